$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Add the four new worksheets at the end of the workbook, in order:
#   Table1, Table2, Table3, "RDS- Postgres DB "
# ---------------------------------------------------------------
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$table1 = $wb.Worksheets.Add($null, $after)
$table1.Name = "Table1"

$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$table2 = $wb.Worksheets.Add($null, $after)
$table2.Name = "Table2"

$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$table3 = $wb.Worksheets.Add($null, $after)
$table3.Name = "Table3"

$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$rds = $wb.Worksheets.Add($null, $after)
$rds.Name = "RDS- Postgres DB "

# ---------------------------------------------------------------
# Table1 : News_Origin / News_URL / RSS_List / Active_Flag
# (cell-fill order matters: it drives the shared-strings table order)
# ---------------------------------------------------------------
$table1.Range("A1").Value = "News_Origin"
$table1.Range("B1").Value = "News_URL"
$table1.Range("D1").Value = "Active_Flag"

$table1.Range("A2").Value = "TOI"
$table1.Range("B2").Value = "https://timesofindia.indiatimes.com/"

$table1.Range("A3").Value = "BS"
$table1.Range("B3").Value = "https://www.business-standard.com/"

$table1.Range("C2").Value = "https://timesofindia.indiatimes.com/rss.cms"
$table1.Range("C1").Value = "RSS_List"
$table1.Range("C3").Value = "https://www.business-standard.com/rss-feeds/listing/"
$table1.Range("C4").Value = "https://www.hindustantimes.com/rss"
$table1.Range("B4").Value = "https://www.hindustantimes.com/"
$table1.Range("A4").Value = "HT"

$table1.Range("A5").Value = "NDTV"
$table1.Range("C5").Value = "https://www.ndtv.com/rss"
$table1.Range("B5").Value = "https://www.ndtv.com/"

$table1.Range("D2").Value = 1
$table1.Range("D3").Value = 0
$table1.Range("D4").Value = 0
$table1.Range("D5").Value = 0

# Real hyperlinks (B2,B3,B4,C4,B5,C5) -- these set both the relationship
# and the "Hyperlink" cell style.
$table1.Hyperlinks.Add($table1.Range("B2"), "https://timesofindia.indiatimes.com/") | Out-Null
$table1.Hyperlinks.Add($table1.Range("B3"), "https://www.business-standard.com/") | Out-Null
$table1.Hyperlinks.Add($table1.Range("B4"), "https://www.hindustantimes.com/") | Out-Null
$table1.Hyperlinks.Add($table1.Range("C4"), "https://www.hindustantimes.com/rss") | Out-Null
$table1.Hyperlinks.Add($table1.Range("B5"), "https://www.ndtv.com/") | Out-Null
$table1.Hyperlinks.Add($table1.Range("C5"), "https://www.ndtv.com/rss") | Out-Null

# C2 / C3 only carry the visual "Hyperlink" style, no actual link.
$table1.Range("C2").Style = "Hyperlink"
$table1.Range("C3").Style = "Hyperlink"

$table1.Columns.Item(1).ColumnWidth = 11.33
$table1.Columns.Item(2).ColumnWidth = 32.11
$table1.Columns.Item(3).ColumnWidth = 46.11
$table1.Columns.Item(4).ColumnWidth = 10.22

$table1.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------
# Table3 : RSS_Feed / Title / Description / Body / Publish_Date /
#          Inserted_Time_Stamp / Sentiment
# (filled first -- it's the source of the new shared strings
#  RSS_Feed/Title/Description/Body/Publish_Date/Sentiment; F1 is
#  filled in last, after a detour through Table2)
# ---------------------------------------------------------------
$table3.Range("A1").Value = "RSS_Feed"
$table3.Range("B1").Value = "Title"
$table3.Range("C1").Value = "Description"
$table3.Range("D1").Value = "Body"
$table3.Range("E1").Value = "Publish_Date"
$table3.Range("G1").Value = "Sentiment"

# ---------------------------------------------------------------
# Table2 : News_URL / RSS_List / RSS_Feed / Active_Flag / Updated_Time_Stamp
# ---------------------------------------------------------------
$table2.Range("A1").Value = "News_URL"
$table2.Range("B1").Value = "RSS_List"
$table2.Range("C1").Value = "RSS_Feed"
$table2.Range("D1").Value = "Active_Flag"
$table2.Range("E1").Value = "Updated_Time_Stamp"

$table3.Range("F1").Value = "Inserted_Time_Stamp"
$table3.Range("J19").Select() | Out-Null

$table2.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------
# RDS- Postgres DB : AWS Postgres RDS connection details
# ---------------------------------------------------------------
$rds.Range("A1").Value = "AWS postgres DB Details:-"
$rds.Range("B2").Value = "DB identifier: postgres-db-identifier1"
$rds.Range("B3").Value = "master user name: postgres_admin"
$rds.Range("B4").Value = "password: postgres123 "
$rds.Range("B5").Value = "Database port : 5432"
$rds.Range("B6").Value = "KMS key ID"
$rds.Range("B7").Value = "d5f76173-6e42-4fae-94fd-e6b8bd9bdc7d"
$rds.Range("B9").Value = "DATABASE NAME: db_news_feed"

$rds.Range("F13").Select() | Out-Null

# ---------------------------------------------------------------
# Make "RDS- Postgres DB " the active (visible) tab, matching the
# commit's bookViews/activeTab + tabSelected state.
# ---------------------------------------------------------------
$rds.Activate()
